# Add new task-report rows 143-158 to the "SB Squares Tasks" sheet.
# (Rows 108 and 121 are pre-existing fully-empty placeholder rows with no
#  cell content; they are naturally dropped from sheetData on save and
#  need no explicit action here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A,B,C,D,E,F ($null = leave blank)
$newRows = @(
    @(143, 20,  "feature", "Floating bottom bar for tentative pick UI (PickControls manual mode redesign, info/warning toasts, grid padding)", "ui-dev", "complete", $null),
    @(144, 21,  "bugfix",  "Replace pick toasts with modals (info toast firing repeatedly, 30s warning not prominent enough)", "ui-dev", "complete", $null),
    @(145, 145, "bugfix",  "Fix dot indicator count showing x/10 instead of x/5 - totalPicks was picksUsed+picksRemaining (double counting tentative picks), changed to picksRemaining only", "ui-dev", "Done", $null),
    @(146, 146, "bugfix",  "Fix timeout handler stale closure and over-filling - added handleTimeoutRef for latest closure, re-fetches fresh squares from DB instead of stale state, added Math.max(0,...) safety cap", "ui-dev", "Done", $null),
    @(147, 147, "bugfix",  "Add mutex lock (pickingRef) to handleTentativePick to prevent race condition from fast double-clicks", "ui-dev", "Done", $null),
    @(148, 131, "Enhanced grid cell visibility (solid borders, glow, thicker width for self-picks in dark mode)", "Enhancement", "ui-dev", "Completed", "components/GridCell.tsx"),
    @(149, 132, "Schema migration: add is_tentative and tentative_started_at columns for tentative picks", "Feature", "architect", "Completed", "supabase/schema.sql, supabase/migrations/20250208000000_tentative_picks.sql, lib/game-logic.ts, app/api/live-scores/route.ts"),
    @(150, 133, "Implement tentative pick logic with circular replacement queue and 2-min timer", "Feature", "ui-dev", "Completed", "app/game/[gameId]/page.tsx"),
    @(151, 134, "Update Grid and GridCell for tentative pick visual states (pulsing ring, pick number badges)", "Feature", "ui-dev", "Completed", "components/Grid.tsx, components/GridCell.tsx"),
    @(152, 135, "Update PickControls with confirm button, countdown timer, and floating bottom bar", "Feature", "ui-dev", "Completed", "components/PickControls.tsx"),
    @(153, 136, "Update admin page for tentative pick cleanup (pickOnBehalf, resetGame, clearPlayerPicks)", "Feature", "architect", "Completed", "app/game/[gameId]/admin/page.tsx"),
    @(154, 137, "Debug tentative picks failure - root cause: migration not applied to deployed DB", "Bug Fix", "architect + ui-dev", "Completed", "N/A (migration deployment)"),
    @(155, 138, "Fix maxPicks derivation using getDraftConfig instead of stale picks_remaining", "Bug Fix", "ui-dev", "Completed", "app/game/[gameId]/page.tsx, components/PickControls.tsx"),
    @(156, 139, "Build admin square reassign feature with mini-grid and player dropdown", "Feature", "architect", "Completed", "app/game/[gameId]/admin/page.tsx"),
    @(157, 140, "Add ref-based mutex to handleRandomPick to prevent double-fire race condition", "Bug Fix", "ui-dev", "Completed", "components/PickControls.tsx"),
    @(158, 158, "bugfix", "Add ref-based mutex to handleRandomPick in PickControls to prevent double-fire via devtools bypass", "ui-dev", "Done", $null)
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    for ($col = 1; $col -le 6; $col++) {
        $val = $entry[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}
